$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 461, shifting existing rows 461:510 down to 462:511
$ws.Rows.Item(461).Insert()

# Populate the new row 461 with its data (mirrors the surrounding rows for the
# constant columns, and carries the new data point's own values)
$ws.Cells.Item(461, 1).Value = 10
$ws.Cells.Item(461, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(461, 3).Value = "La Araucanía"
$ws.Cells.Item(461, 4).Value = 44918
$ws.Cells.Item(461, 5).Value = 9
$ws.Cells.Item(461, 6).Value = 100114014
$ws.Cells.Item(461, 7).Value = "Betarraga"
$ws.Cells.Item(461, 8).Value = "Sin especificar"
$ws.Cells.Item(461, 9).Value = "Primera"
$ws.Cells.Item(461, 10).Value = 95
$ws.Cells.Item(461, 11).Value = 10000
$ws.Cells.Item(461, 12).Value = 10000
$ws.Cells.Item(461, 13).Value = 10000
$ws.Cells.Item(461, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(461, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(461, 16).Value = 833
$ws.Cells.Item(461, 17).Value = 12
$ws.Cells.Item(461, 18).Value = "Hortaliza"
